$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1619.339
$ws.Range("J112").Value = 1621.8276
$ws.Range("L112").Value = 4865.4828
$ws.Range("N112").Value = -7081.4828
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11486.434
$ws.Range("I2").Value = 12355.296
$ws.Range("K2").Value = 12355.296
$ws.Range("M2").Value = -12242.296
$ws.Range("H32").Value = 10933.088
$ws.Range("I32").Value = 7866
$ws.Range("K32").Value = 7866
$ws.Range("M32").Value = -7579
$ws.Range("H61").Value = 4620.9565
$ws.Range("I61").Value = 3363.4
$ws.Range("K61").Value = 3363.4
$ws.Range("M61").Value = -3151.4
$ws.Range("H63").Value = 4277.619
$ws.Range("I63").Value = 4175.4736
$ws.Range("K63").Value = 4175.4736
$ws.Range("M63").Value = -3489.4736
$ws.Range("H66").Value = 4277.619
$ws.Range("I66").Value = 4175.4736
$ws.Range("K66").Value = 20877.368
$ws.Range("M66").Value = -17445.368
$ws.Range("H74").Value = 2260.8147
$ws.Range("J74").Value = 4006.75
$ws.Range("L74").Value = 4006.75
$ws.Range("N74").Value = -5754.75
$ws.Range("H77").Value = 2260.8147
$ws.Range("J77").Value = 4006.75
$ws.Range("L77").Value = 20033.75
$ws.Range("N77").Value = -28769.75
$ws.Range("H97").Value = 2776.75
$ws.Range("I97").Value = 2702.3333
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 2702.3333
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
$ws.Range("M97").Value = -2206.3333
$ws.Range("H116").Value = 11486.434
$ws.Range("I116").Value = 12355.296
$ws.Range("K116").Value = 12355.296
$ws.Range("M116").Value = -10061.296
$ws.Range("H122").Value = 4818.469
$ws.Range("I122").Value = 4040.6775
$ws.Range("J122").Value = 6158
$ws.Range("K122").Value = 12122.0325
$ws.Range("L122").Value = 18474
$ws.Range("M122").Value = -9672.0325
$ws.Range("N122").Value = -23374
$ws.Range("H136").Value = 4620.9565
$ws.Range("I136").Value = 3363.4
$ws.Range("K136").Value = 10090.2
$ws.Range("M136").Value = -7540.200000000001
$ws.Range("H138").Value = 86666.664
$ws.Range("J138").Value = 86666.664
$ws.Range("L138").Value = 86666.664
$ws.Range("N138").Value = -96946.664
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11486.434
$ws.Range("I3").Value = 12355.296
$ws.Range("K3").Value = 12355.296
$ws.Range("M3").Value = -12241.296
$ws.Range("H97").Value = 9888.571
$ws.Range("I97").Value = 9888.571
$ws.Range("K97").Value = 9888.571
$ws.Range("M97").Value = -8897.571
$ws.Range("H105").Value = 2870.0625
$ws.Range("I105").Value = 1829.1818
$ws.Range("J105").Value = 5160
$ws.Range("K105").Value = 1829.1818
$ws.Range("L105").Value = 5160
$ws.Range("M105").Value = -82.18180000000007
$ws.Range("N105").Value = -8654
$ws.Range("H134").Value = 4138.241
$ws.Range("I134").Value = 2683.0557
$ws.Range("J134").Value = 6519.4546
$ws.Range("K134").Value = 8049.1671
$ws.Range("L134").Value = 19558.3638
$ws.Range("M134").Value = -5514.1671
$ws.Range("N134").Value = -24628.3638
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3772.1667
$ws.Range("I31").Value = 919.5625
$ws.Range("K31").Value = 919.5625
$ws.Range("M31").Value = -624.5625
$ws.Range("H34").Value = 3772.1667
$ws.Range("I34").Value = 919.5625
$ws.Range("K34").Value = 919.5625
$ws.Range("M34").Value = -717.5625
$ws.Range("H58").Value = 3988.2727
$ws.Range("I58").Value = 2974.75
$ws.Range("J58").Value = 4567.4287
$ws.Range("K58").Value = 2974.75
$ws.Range("L58").Value = 4567.4287
$ws.Range("M58").Value = -2771.75
$ws.Range("N58").Value = -4973.4287
$ws.Range("I105").Value = 22222692
$ws.Range("K105").Value = 22222692
$ws.Range("M105").Value = -22220945
$ws.Range("H136").Value = 3988.2727
$ws.Range("I136").Value = 2974.75
$ws.Range("J136").Value = 4567.4287
$ws.Range("K136").Value = 8924.25
$ws.Range("L136").Value = 13702.2861
$ws.Range("M136").Value = -6374.25
$ws.Range("N136").Value = -18802.2861
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 116
$ws.Range("J33").Value = 85
$ws.Range("L33").Value = 510
$ws.Range("N33").Value = -1076
$ws.Range("H39").Value = 3290.75
$ws.Range("J39").Value = 3483.4666
$ws.Range("L39").Value = 10450.3998
$ws.Range("N39").Value = -11038.3998
$ws.Range("H55").Value = 717822
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 717822
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 2153466
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -2153820
$ws.Range("H107").Value = 14493031
$ws.Range("I107").Value = 356
$ws.Range("J107").Value = 16666932
$ws.Range("K107").Value = 1068
$ws.Range("L107").Value = 50000796
$ws.Range("M107").Value = 852
$ws.Range("N107").Value = -50004636
$ws.Range("H131").Value = 5973.75
$ws.Range("I131").Value = 3756.25
$ws.Range("J131").Value = 8191.25
$ws.Range("K131").Value = 11268.75
$ws.Range("L131").Value = 24573.75
$ws.Range("M131").Value = -6228.75
$ws.Range("N131").Value = -34653.75
$ws.Range("H132").Value = 699.75
$ws.Range("I132").Value = 599.5
$ws.Range("K132").Value = 5395.5
$ws.Range("M132").Value = -2865.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 26996
$ws.Range("H80").Value = 52594.137
$ws.Range("I80").Value = 71030.06
$ws.Range("K80").Value = 71030.06
$ws.Range("M80").Value = -70032.06
$ws.Range("H83").Value = 52594.137
$ws.Range("I83").Value = 71030.06
$ws.Range("K83").Value = 355150.3
$ws.Range("M83").Value = -350158.3
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 458203.22
$ws.Range("I93").Value = 3518.2
$ws.Range("J93").Value = 1432528.2
$ws.Range("K93").Value = 3518.2
$ws.Range("L93").Value = 1432528.2
$ws.Range("M93").Value = -2270.2
$ws.Range("N93").Value = -1435024.2
$ws.Range("H132").Value = 3693.2368
$ws.Range("I132").Value = 3024.926
$ws.Range("J132").Value = 5333.636
$ws.Range("K132").Value = 9074.778
$ws.Range("L132").Value = 16000.908
$ws.Range("M132").Value = -6544.778
$ws.Range("N132").Value = -21060.908
$ws.Range("H136").Value = 4986.6284
$ws.Range("I136").Value = 3774.1
$ws.Range("K136").Value = 11322.3
$ws.Range("M136").Value = -8772.3
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 45164.668
$ws.Range("J49").Value = 45164.668
$ws.Range("L49").Value = 45164.668
$ws.Range("N49").Value = -45624.668
$ws.Range("H100").Value = 698.61536
$ws.Range("I100").Value = 698.61536
$ws.Range("K100").Value = 1397.23072
$ws.Range("M100").Value = -856.23072
$ws.Range("H132").Value = 1917.5555
$ws.Range("I132").Value = 1359.3489
$ws.Range("K132").Value = 4078.0467
$ws.Range("M132").Value = -1548.0467
$ws.Range("H136").Value = 4060.6155
$ws.Range("I136").Value = 3394.125
$ws.Range("K136").Value = 10182.375
$ws.Range("M136").Value = -7632.375
